# Commit: "Change Excel Field View to Cache, And set default value to FALSE"
#
# The "Property" sheet has a header row (row 1) describing field names for
# each data column, and column F is the "View" flag, currently defaulted to
# TRUE (1) for every data row. This edit:
#   1. Renames the column F header from "View" to "Cache".
#   2. Resets the default value of column F (rows 2-15) from TRUE to FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header F1: "View" -> "Cache"
$ws.Range("F1").Value = "Cache"

# 2) Flip the default for every data row in column F from TRUE to FALSE
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
}
